$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D3 was "Complete" -> now "Suited to Manual"
$ws.Range("D3").Value = "Suited to Manual"

# E3 used to hold "Suite to manual" -> that information moved into D3,
# so the cell is now fully cleared out (content + formatting).
$ws.Range("E3").Clear()

# Update the active selection on the sheet to E3.
$ws.Activate()
$ws.Range("E3").Select()
